$wb = $excel.ActiveWorkbook

# --- Work on the "Repayment Schedule" sheet (3rd sheet) ---
$ws = $wb.Worksheets.Item(3)

# Insert a new (blank) column before column N. This shifts the old
# "Late" column (N) to O, the old (empty) O to P, and the old
# "Outstanding" column (P) to Q - matching the new header layout:
# ... M=In Advance, N=(blank), O=Late, P=Heading, Q=Outstanding
$ws.Columns.Item(14).Insert()

# Activate the "Repayment Schedule" sheet and select M17, matching the
# saved selection/active-tab state from the authored workbook.
$ws.Activate()
$ws.Range("M17").Select()
